$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (date D2 and volume J2 change)
$ws.Range("D2").Value = 44848
$ws.Range("J2").Value = 1000

# Row 4 (date D4, K4, M4, P4 change)
$ws.Range("D4").Value = 44881
$ws.Range("K4").Value = 1900
$ws.Range("M4").Value = 1950
$ws.Range("P4").Value = 650

# Row 5 (date D5, J5, K5, M5, P5 change)
$ws.Range("D5").Value = 44685
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 1500
$ws.Range("M5").Value = 1750
$ws.Range("P5").Value = 583

# Row 6 (date D6, J6, K6, M6, P6 change)
$ws.Range("D6").Value = 44883
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 1800
$ws.Range("M6").Value = 1900
$ws.Range("P6").Value = 633
